$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 12; this shifts the existing rows 12..113
# down to 13..114 (matching the diff's row-shift pattern) and grows the
# sheet dimension to A1:R114 automatically.
$ws.Rows("12:12").Insert()

# Populate the newly inserted row 12 with this week's data.
$ws.Range("A12").Value = 5
$ws.Range("B12").Value = "Macroferia Regional de Talca"
$ws.Range("C12").Value = "Maule"
$ws.Range("D12").Value = 44530
$ws.Range("E12").Value = 7
$ws.Range("F12").Value = 100112031
$ws.Range("G12").Value = "Poroto verde"
$ws.Range("H12").Value = "Sin especificar"
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 150
$ws.Range("K12").Value = 20000
$ws.Range("L12").Value = 20000
$ws.Range("M12").Value = 20000
$ws.Range("N12").Value = "$/saco 25 kilos"
$ws.Range("O12").Value = "Región del Maule"
$ws.Range("P12").Value = 800
$ws.Range("Q12").Value = 25
$ws.Range("R12").Value = "Hortaliza"
